$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 667, shifting existing rows 667:744 down to 668:745.
$ws.Rows(667).Insert()

# Populate the newly inserted row 667 with the new record's data.
$ws.Range("A667").Value = 3
$ws.Range("B667").Value = "Femacal de La Calera"
$ws.Range("C667").Value = "Coquimbo"
$ws.Range("D667").Value = 45142
$ws.Range("E667").Value = 5
$ws.Range("F667").Value = 100112021
$ws.Range("G667").Value = "Ají"
$ws.Range("H667").Value = "Inferno"
$ws.Range("I667").Value = "Primera"
$ws.Range("J667").Value = 40
$ws.Range("K667").Value = 14000
$ws.Range("L667").Value = 14000
$ws.Range("M667").Value = 14000
$ws.Range("N667").Value = "$/caja 10 kilos"
$ws.Range("O667").Value = "Región de Arica y Parinacota"
$ws.Range("P667").Value = 1400
$ws.Range("Q667").Value = 10
$ws.Range("R667").Value = "Hortaliza"
